$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '310.79'
$ws.Range("E2").Value = '-1.39%'
$ws.Range("D3").Value = '37.54'
$ws.Range("E3").Value = '-4.34%'
$ws.Range("D4").Value = '5.087'
$ws.Range("E4").Value = '-1.15%'
$ws.Range("D5").Value = '0.07752'
$ws.Range("E5").Value = '-5.01%'
$ws.Range("D6").Value = '4.346'
$ws.Range("E6").Value = '-1.42%'
$ws.Range("E7").Value = '-1.58%'
$ws.Range("E8").Value = '-4.69%'
$ws.Range("E9").Value = '-11.69%'
$ws.Range("D10").Value = '0.9198'
$ws.Range("E10").Value = '-1.74%'
$ws.Range("D11").Value = '0.1190'
$ws.Range("E11").Value = '-9.30%'
$ws.Range("D12").Value = '0.1918'
$ws.Range("E12").Value = '-3.44%'
$ws.Range("D13").Value = '0.08879'
$ws.Range("E13").Value = '-1.48%'
$ws.Range("D14").Value = '0.03387'
$ws.Range("E14").Value = '-2.93%'
$ws.Range("D15").Value = '0.09703'
$ws.Range("E15").Value = '-0.19%'
$ws.Range("D16").Value = '0.001369'
$ws.Range("E16").Value = '-2.86%'
$ws.Range("D17").Value = '0.005858'
$ws.Range("E17").Value = '-1.92%'
$ws.Range("D18").Value = '3.554'
$ws.Range("E18").Value = '-1.51%'
$ws.Range("E19").Value = '-1.77%'
$ws.Range("D20").Value = '5.036'
$ws.Range("E20").Value = '0.17%'
$ws.Range("D21").Value = '0.1261'
$ws.Range("E21").Value = '-3.66%'
$ws.Range("D22").Value = '0.2591'
$ws.Range("E22").Value = '3.96%'
$ws.Range("D23").Value = '0.02104'
$ws.Range("E23").Value = '5,587.55%'
$ws.Range("D24").Value = '0.04387'
$ws.Range("E24").Value = '0.35%'
$ws.Range("E25").Value = '-2.71%'
$ws.Range("D26").Value = '0.004243'
$ws.Range("E26").Value = '-10.63%'
$ws.Range("E27").Value = '-65.35%'
$ws.Range("D39").Value = '0.02108'
$ws.Range("E39").Value = '-5.68%'
$ws.Range("D40").Value = '0.04948'
$ws.Range("E40").Value = '-5.42%'
$ws.Range("D41").Value = '0.007661'
$ws.Range("E41").Value = '-0.27%'
$ws.Range("D42").Value = '0.009910'
$ws.Range("E42").Value = '-4.11%'
$ws.Range("D43").Value = '0.1342'
$ws.Range("E43").Value = '-3.80%'
$ws.Range("D44").Value = '0.002061'
$ws.Range("E44").Value = '-2.05%'
$ws.Range("D45").Value = '0.009623'
$ws.Range("E45").Value = '5.48%'
$ws.Range("D46").Value = '0.00006569'
$ws.Range("E46").Value = '-3.79%'
$ws.Range("E47").Value = '-0.15%'
$ws.Range("E48").Value = '1.08%'
$ws.Range("E50").Value = '-0.15%'
$ws.Range("E51").Value = '-0.15%'
